$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1
# header cell so the new column matches the other headers' style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Add the new "Save" header and its data column.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
